$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "Username" column (column A), shifting everything left
$ws.Columns.Item(1).Delete()

# Update row 2 (now occupies what used to be columns B-H) with new sample data
$ws.Range("C2").Value = "r"
$ws.Range("A2").Value = "EEEE"
$ws.Range("B2").Value = "FFFF"

# Update header row (row 1) - add Firstname/LastName columns, reorder Password column
$ws.Range("A1").Value = "Firstname"
$ws.Range("B1").Value = "LastName"
$ws.Range("E1").Value = "Password"

$ws.Range("E2").Value = "Hs622!@ad"
$ws.Range("F2").Value = "Hs622!@ad"

# Rebuild hyperlinks: the old ones need to move from column D to column C (after
# the column delete) and two brand-new ones are needed on row 2
$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("C2"), "mailto:DUMMY@EMAIL", [Type]::Missing, "DUMMY@EMAIL")
$ws.Hyperlinks.Add($ws.Range("C4"), "mailto:DUMMY@EMAIL")
$ws.Hyperlinks.Add($ws.Range("C6"), "mailto:DUMMY@EMAIL")
$ws.Hyperlinks.Add($ws.Range("C8"), "mailto:DUMMY@EMAIL")
$ws.Hyperlinks.Add($ws.Range("E2"), "mailto:DUMMY@EMAIL")
$ws.Hyperlinks.Add($ws.Range("F2"), "mailto:DUMMY@EMAIL")

# Restore the Hyperlink cell style (Add() above creates a duplicate style record;
# re-applying the named style keeps every linked cell on the original style index)
$ws.Range("C2").Style = "Hyperlink"
$ws.Range("C4").Style = "Hyperlink"
$ws.Range("C6").Style = "Hyperlink"
$ws.Range("C8").Style = "Hyperlink"
$ws.Range("E2").Style = "Hyperlink"
$ws.Range("F2").Style = "Hyperlink"

$ws.Range("E3").Select()
